$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.978.39"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.646.88"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.19"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.526"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.65"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.877.89"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "1.650.78"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.74"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "27.951.46"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.22"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.72"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.40"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.74"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D33").Value = "1.456.98"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.890"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.920"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.52"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.23"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  +5.76%  "
$ws.Range("D48").Value = "1.786.70"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.78"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("E50").Value = "  +2.05%  "
$ws.Range("E51").Value = "  +1.33%  "
